# Revert "inputs changed from fastqs to bam": WES assay run template.
# The WES sheet's Samples section previously took a single "Bam" file path
# column; this restores the original paired-fastq layout (Forward fastq /
# Reverse fastq) and the matching Legend / comment text.

$wb = $excel.ActiveWorkbook
$wes = $wb.Worksheets.Item("WES")
$legend = $wb.Worksheets.Item("Legend")

# ---------------------------------------------------------------------------
# 1. WES sheet: insert a new column D ("Reverse fastq") before the old
#    "Sequencing date" column, so the Samples table becomes:
#    Cimac id | Forward fastq | Reverse fastq | Sequencing date | Quality flag
# ---------------------------------------------------------------------------
$wes.Columns.Item(4).Insert()

# Rename old "Bam" header (col C) to "Forward fastq", add the new
# "Reverse fastq" header in col D (now vacated by the insert).
$wes.Range("C12").Value = "Forward fastq"
$wes.Range("D12").Value = "Reverse fastq"

# Sample data rows: restore the two fastq paths per sample (previously one
# bam path each).
$wes.Range("C13").Value = "/local/path/to/fwd.1.1.1.fastq.gz"
$wes.Range("D13").Value = "/local/path/to/rev.1.1.1.fastq.gz"
$wes.Range("C14").Value = "/local/path/to/fwd.1.2.1.fastq.gz"
$wes.Range("D14").Value = "/local/path/to/rev.1.2.1.fastq.gz"

# The new D column needs the same "path to a file" comment the forward
# fastq column (C) already carries.
$wes.Range("D12").AddComment("Path to a file on a user's computer.")

# Selection / scroll position, matching the reverted template.
$wes.Activate()
$wes.Range("D15").Select()
$wes.Application.ActiveWindow.ScrollRow = 2

# ---------------------------------------------------------------------------
# 2. Legend sheet: add back the "Assay run id" definition row, and split the
#    single "Bam" row into "Forward fastq" / "Reverse fastq" rows.
# ---------------------------------------------------------------------------
$legend.Rows.Item(3).Insert()
$legend.Range("B3").Value = "Assay run id"
$legend.Range("C3").Value = "String"
$legend.Range("D3").Value = "User defined unique identifier for this assay run."

# After the row-3 insert, the old "Bam" legend row (previously row 13) is
# now row 14. Split it into two rows: "Forward fastq" (row 14) and the new
# "Reverse fastq" (row 15).
$legend.Rows.Item(15).Insert()
$legend.Range("B14").Value = "Forward fastq"
$legend.Range("C14").Value = "String"
$legend.Range("D14").Value = "Path to a file on a user's computer."
$legend.Range("B15").Value = "Reverse fastq"
$legend.Range("B15").Font.Bold = $true
$legend.Range("C15").Value = "String"
$legend.Range("D15").Value = "Path to a file on a user's computer."
